# Applies the "Updated to Old-New Response version (still testing, not final)" edit:
#  - Sheet "Sound_data": scroll view so BH column area is visible, select whole column BH
#  - Sheet "Sheet1": clear the old scratch topLeftCell, select K52, populate column G
#    (rows 1-15, 17-30, 32-44, 45-52) with copies of column D's entries - grouped by their
#    fill-color style (a "sort column D by fill color" paste into column G) - and remove the
#    now-stale raw number block that used to live in G50:J53.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 (tab name "Sound_data") - just a view/selection change
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sound_data")
$ws1.Activate()
$ws1.Range("BH1:BH1048576").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 54

# ---------------------------------------------------------------------------
# Sheet2 (tab name "Sheet1") - the main data edit
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Activate()

# Drop the old leftover numeric scratch block in G50:J53 before writing the
# new column G content (new values land on top of G50:G52, G53 stays empty).
$ws2.Range("G50:J53").ClearContents()

# Column G is filled with copies of column D's cells, grouped by the D
# cell's fill-color style, preserving top-to-bottom order within each group.
# Group 1 (style: light-blue-ish) -> G1:G15
$ws2.Range("D4").Copy($ws2.Range("G1"))
$ws2.Range("D7").Copy($ws2.Range("G2"))
$ws2.Range("D8").Copy($ws2.Range("G3"))
$ws2.Range("D9").Copy($ws2.Range("G4"))
$ws2.Range("D11").Copy($ws2.Range("G5"))
$ws2.Range("D12").Copy($ws2.Range("G6"))
$ws2.Range("D14").Copy($ws2.Range("G7"))
$ws2.Range("D22").Copy($ws2.Range("G8"))
$ws2.Range("D29").Copy($ws2.Range("G9"))
$ws2.Range("D35").Copy($ws2.Range("G10"))
$ws2.Range("D37").Copy($ws2.Range("G11"))
$ws2.Range("D43").Copy($ws2.Range("G12"))
$ws2.Range("D46").Copy($ws2.Range("G13"))
$ws2.Range("D52").Copy($ws2.Range("G14"))
$ws2.Range("D54").Copy($ws2.Range("G15"))

# Group 2 -> G17:G30 (row 16 left blank as a spacer)
$ws2.Range("D1").Copy($ws2.Range("G17"))
$ws2.Range("D2").Copy($ws2.Range("G18"))
$ws2.Range("D3").Copy($ws2.Range("G19"))
$ws2.Range("D6").Copy($ws2.Range("G20"))
$ws2.Range("D13").Copy($ws2.Range("G21"))
$ws2.Range("D20").Copy($ws2.Range("G22"))
$ws2.Range("D24").Copy($ws2.Range("G23"))
$ws2.Range("D25").Copy($ws2.Range("G24"))
$ws2.Range("D27").Copy($ws2.Range("G25"))
$ws2.Range("D32").Copy($ws2.Range("G26"))
$ws2.Range("D33").Copy($ws2.Range("G27"))
$ws2.Range("D34").Copy($ws2.Range("G28"))
$ws2.Range("D48").Copy($ws2.Range("G29"))
$ws2.Range("D50").Copy($ws2.Range("G30"))

# Group 3 -> G32:G44 (row 31 left blank as a spacer)
$ws2.Range("D5").Copy($ws2.Range("G32"))
$ws2.Range("D10").Copy($ws2.Range("G33"))
$ws2.Range("D15").Copy($ws2.Range("G34"))
$ws2.Range("D16").Copy($ws2.Range("G35"))
$ws2.Range("D17").Copy($ws2.Range("G36"))
$ws2.Range("D18").Copy($ws2.Range("G37"))
$ws2.Range("D19").Copy($ws2.Range("G38"))
$ws2.Range("D21").Copy($ws2.Range("G39"))
$ws2.Range("D23").Copy($ws2.Range("G40"))
$ws2.Range("D26").Copy($ws2.Range("G41"))
$ws2.Range("D28").Copy($ws2.Range("G42"))
$ws2.Range("D30").Copy($ws2.Range("G43"))
$ws2.Range("D31").Copy($ws2.Range("G44"))

# Group 4 -> G45:G46, G48:G52, with the last item of group 3 trailing at G47
$ws2.Range("D36").Copy($ws2.Range("G45"))
$ws2.Range("D38").Copy($ws2.Range("G46"))
$ws2.Range("D39").Copy($ws2.Range("G47"))
$ws2.Range("D40").Copy($ws2.Range("G48"))
$ws2.Range("D42").Copy($ws2.Range("G49"))
$ws2.Range("D44").Copy($ws2.Range("G50"))
$ws2.Range("D45").Copy($ws2.Range("G51"))
$ws2.Range("D47").Copy($ws2.Range("G52"))

# Final selection / scroll position on this sheet
$ws2.Range("K52").Select()

# Leave "Sheet1" (sheet2) as the active tab, matching the source file.
$ws2.Activate()
